$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pedido 69134d11b9c1d30b15fabdc3:
# The original "Test Ringover (NO TOCAR)" line (row 3) had some leftover
# empty placeholder cells (Unidades Estructura/Paneles, Optimizador,
# Unidades Optimizador, Cargador VE) that are no longer needed - clear them.
$ws.Cells.Item(3, 5).ClearContents()   # E3 - Unidades Estructura/Paneles
$ws.Cells.Item(3, 6).ClearContents()   # F3 - Optimizador
$ws.Cells.Item(3, 7).ClearContents()   # G3 - Unidades Optimizador
$ws.Cells.Item(3, 12).ClearContents()  # L3 - Cargador VE

# Add a new order line (row 4), duplicating row 3's data but referencing the
# BATERÍA LITIO SIGEN ENERGY SIGENSTOR 10,0KW battery (qty 3) instead of the
# 8,0KW battery (qty 1).
$ws.Cells.Item(4, 1).Value = 2488
$ws.Cells.Item(4, 2).Value = "Test Ringover (NO TOCAR)"
$ws.Cells.Item(4, 3).Value = "Estructura coplanar NOVOTEGRA"

# Columns D, I and K hold numeric-looking quantities that must stay stored
# as text (like the rest of the sheet), so force text formatting before
# assigning and then drop the formatting again so no style gets attached.
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1"
$ws.Cells.Item(4, 4).ClearFormats()

$ws.Cells.Item(4, 8).Value = "Inversor híbrido monofásico SUN-6k-SG05LP1-EU"

$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = "1"
$ws.Cells.Item(4, 9).ClearFormats()

$ws.Cells.Item(4, 10).Value = "BATERÍA LITIO SIGEN ENERGY SIGENSTOR 10,0KW"

$ws.Cells.Item(4, 11).NumberFormat = "@"
$ws.Cells.Item(4, 11).Value = "3"
$ws.Cells.Item(4, 11).ClearFormats()

$ws.Cells.Item(4, 13).Value = "Sí"
$ws.Cells.Item(4, 14).Value = "2024-01-03T10:49:29.104Z"
